# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column right before "GZ" (column 208) so that the
# existing "nom" (GZ) / "url_produit" (HA) columns shift one place to the
# right, becoming HA / HB respectively. This mirrors a new price-check
# column being added to the price-history table, just before the two
# trailing metadata columns.
$ws.Columns.Item(208).Insert()

# Header row: the freshly inserted GZ1 cell gets the new check timestamp
# (the same bold/centered header style carries over automatically from
# the Insert shift).
$ws.Range("GZ1").Value = "2026-02-06 15:29:59"

# Data rows: populate the new GZ column with the same price value already
# recorded in the previous check column (GY) - i.e. the price observed at
# the prior check is carried forward unchanged into this new check column.
# Rows where GY has no recorded price are left blank, matching GY.
$lastRow = 210
for ($r = 2; $r -le $lastRow; $r++) {
    $prev = $ws.Range("GY$r").Value2
    if ($prev -ne $null -and $prev -ne "") {
        $ws.Range("GZ$r").Value = $prev
    }
}
